$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A80").Value = "GRT-USD"
